# Adding Astro and IMperial to the Portfolio Forecast
# Refresh the forecasted-consumption series (column A) with the latest
# model run and advance the 15-minute timestamp series (column B) to the
# new forecast day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$forecast = @(5430, 5380, 5330, 5270, 5230, 5210, 5190, 5170, 5150, 5150, 5150, 5150, 5170, 5170, 5190, 5210, 5250, 5300, 5350, 5420, 5500, 5610, 5730, 5850, 6050, 6150, 6230, 6280, 6290, 6270, 6220, 6140, 5970, 5830, 5690, 5540, 5420, 5280, 5160, 5050, 4910, 4830, 4760, 4710, 4670, 4650, 4620, 4600, 4570, 4550, 4550, 4550, 4550, 4570, 4600, 4620, 4630, 4670, 4730, 4780, 4890, 4990, 5100, 5220, 5320, 5450, 5560, 5670, 5820, 5930, 6050, 6190, 6290, 6440, 6580, 6690, 6830, 6950, 7060, 7160, 7250, 7240, 7210, 7100, 6880, 6730, 6590, 6430, 6200, 6020, 5930, 5820, 5530, 5500, 5430, 5350)

$timestamps = @(45776, 45776.01041666666, 45776.02083333334, 45776.03125, 45776.04166666666, 45776.05208333334, 45776.0625, 45776.07291666666, 45776.08333333334, 45776.09375, 45776.10416666666, 45776.11458333334, 45776.125, 45776.13541666666, 45776.14583333334, 45776.15625, 45776.16666666666, 45776.17708333334, 45776.1875, 45776.19791666666, 45776.20833333334, 45776.21875, 45776.22916666666, 45776.23958333334, 45776.25, 45776.26041666666, 45776.27083333334, 45776.28125, 45776.29166666666, 45776.30208333334, 45776.3125, 45776.32291666666, 45776.33333333334, 45776.34375, 45776.35416666666, 45776.36458333334, 45776.375, 45776.38541666666, 45776.39583333334, 45776.40625, 45776.41666666666, 45776.42708333334, 45776.4375, 45776.44791666666, 45776.45833333334, 45776.46875, 45776.47916666666, 45776.48958333334, 45776.5, 45776.51041666666, 45776.52083333334, 45776.53125, 45776.54166666666, 45776.55208333334, 45776.5625, 45776.57291666666, 45776.58333333334, 45776.59375, 45776.60416666666, 45776.61458333334, 45776.625, 45776.63541666666, 45776.64583333334, 45776.65625, 45776.66666666666, 45776.67708333334, 45776.6875, 45776.69791666666, 45776.70833333334, 45776.71875, 45776.72916666666, 45776.73958333334, 45776.75, 45776.76041666666, 45776.77083333334, 45776.78125, 45776.79166666666, 45776.80208333334, 45776.8125, 45776.82291666666, 45776.83333333334, 45776.84375, 45776.85416666666, 45776.86458333334, 45776.875, 45776.88541666666, 45776.89583333334, 45776.90625, 45776.91666666666, 45776.92708333334, 45776.9375, 45776.94791666666, 45776.95833333334, 45776.96875, 45776.97916666666, 45776.98958333334)

for ($i = 0; $i -lt $forecast.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $forecast[$i]
    $ws.Cells.Item($row, 2).Value = $timestamps[$i]
}
